# Apply the updated cryptocurrency price/volume figures (and the two
# coin re-rankings that swapped adjacent rows) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.117.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.058.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.673'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.35'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.379'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0807'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.27%  '
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.357.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.815'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.056.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.041.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0930'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +12.33%  '
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '74.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  +1.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = '  -3.08%  '
$ws.Range("E33").Value = '  +2.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.86%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0873'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.68%  '
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.107'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +23.39%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.46%  '
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +54.05%  '
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '97.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.39%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.301.10'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.59%  '
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.33%  '
